# set github as mirror
# Insert a duplicate header row at row 11, pushing the rcv1subset3..yeast rows
# down by one (rows 11-18 -> 12-19), and update the named range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (shifts existing rows 11..18 down to 12..19)
$ws.Rows.Item(11).Insert()

# Populate the newly-inserted row 11 with a repeat of the header row (row 1)
$ws.Range("A11").Value = "name"
$ws.Range("B11").Value = "domain"
$ws.Range("C11").Value = "features_num"
$ws.Range("D11").Value = "labels_num"
$ws.Range("E11").Value = "instances_train"
$ws.Range("F11").Value = "instances_test"

# Update the defined name "available_data" to reflect the new extent
$wb.Names.Item("available_data").RefersTo = "=Sheet1!`$A`$1:`$F`$19"

# Update the active selection as per the saved workbook state
$ws.Range("C9").Select()
